$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.370.28"
$ws.Range("E2").Value = "  -2.59%  "

$ws.Range("D3").Value = "3.183.13"
$ws.Range("E3").Value = "  -4.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.28%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.181.77"
$ws.Range("E8").Value = "  -4.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.501"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.28%  "

$ws.Range("E10").Value = "  -5.87%  "

$ws.Range("E11").Value = "  -6.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.42%  "

$ws.Range("E13").Value = "  -6.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.29%  "

$ws.Range("D15").Value = "3.707.74"
$ws.Range("E15").Value = "  -4.16%  "

$ws.Range("E16").Value = "  -1.44%  "

$ws.Range("D17").Value = "3.184.81"
$ws.Range("E17").Value = "  -3.99%  "

$ws.Range("D18").Value = "62.406.76"
$ws.Range("E18").Value = "  -2.68%  "

$ws.Range("E19").Value = "  -5.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "456.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.31%  "

$ws.Range("E22").Value = "  -5.09%  "

$ws.Range("E23").Value = "  -5.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.06%  "

$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.79%  "

$ws.Range("E31").Value = "  -7.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.102"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.89%  "

$ws.Range("E35").Value = "  -6.29%  "

$ws.Range("E36").Value = "  -3.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.04%  "

$ws.Range("E38").Value = "  -10.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0385"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "410.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.09%  "

$ws.Range("D41").Value = "2.941.60"
$ws.Range("E41").Value = "  -4.05%  "

$ws.Range("E42").Value = "  +0.75%  "

$ws.Range("E43").Value = "  -5.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.93%  "

$ws.Range("E45").Value = "  -3.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.249"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "

